$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = 111979876
$ws.Range("B5").Value = 90666
$ws.Range("C5").Value = "Ovaliderad"
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 4364
$ws.Range("F5").Value = "Dropptaggsvamp"
$ws.Range("G5").Value = "Hydnellum ferrugineum"
$ws.Range("H5").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "10"
$ws.Range("P5").Value = "Strömsbergs bruk , Upl"
$ws.Range("Q5").Value = 642175.212860164
$ws.Range("R5").Value = 6698319.320043332
$ws.Range("S5").Value = 25
$ws.Range("T5").Value = "Uppsala"
$ws.Range("U5").Value = "Tierp"
$ws.Range("V5").Value = "Uppland"
$ws.Range("W5").Value = "Tolfta"
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = "2023-09-09"
$ws.Range("Z5").NumberFormat = "@"
$ws.Range("Z5").Value = "11:56"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = "2023-09-09"
$ws.Range("AB5").NumberFormat = "@"
$ws.Range("AB5").Value = "11:56"
$ws.Range("AD5").Value = $False
$ws.Range("AE5").Value = $False
$ws.Range("AG5").Value = $False
$ws.Range("AW5").Value = "Lotta Lund"
$ws.Range("AX5").Value = "Lotta Lund"

# Row 6
$ws.Range("A6").Value = 111979480
$ws.Range("B6").Value = 88909
$ws.Range("C6").Value = "Ovaliderad"
$ws.Range("D6").Value = "VU"
$ws.Range("E6").Value = 720
$ws.Range("F6").Value = "Violgubbe"
$ws.Range("G6").Value = "Gomphus clavatus"
$ws.Range("H6").Value = "(Pers.) Gray"
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "70"
$ws.Range("P6").Value = "Strömsbergs bruk , Upl"
$ws.Range("Q6").Value = 642175.212860164
$ws.Range("R6").Value = 6698319.320043332
$ws.Range("S6").Value = 25
$ws.Range("T6").Value = "Uppsala"
$ws.Range("U6").Value = "Tierp"
$ws.Range("V6").Value = "Uppland"
$ws.Range("W6").Value = "Tolfta"
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = "2023-09-09"
$ws.Range("Z6").NumberFormat = "@"
$ws.Range("Z6").Value = "11:47"
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "2023-09-09"
$ws.Range("AB6").NumberFormat = "@"
$ws.Range("AB6").Value = "11:47"
$ws.Range("AD6").Value = $False
$ws.Range("AE6").Value = $False
$ws.Range("AG6").Value = $False
$ws.Range("AH6").Value = "Skogsmark"
$ws.Range("AW6").Value = "Lotta Lund"
$ws.Range("AX6").Value = "Lotta Lund"

# Row 7
$ws.Range("A7").Value = 111980195
$ws.Range("B7").Value = 90684
$ws.Range("C7").Value = "Ovaliderad"
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 4368
$ws.Range("F7").Value = "Dofttaggsvamp"
$ws.Range("G7").Value = "Hydnellum suaveolens"
$ws.Range("H7").Value = "(Scop.:Fr.) P. Karst."
$ws.Range("P7").Value = "Tierp, Upl"
$ws.Range("Q7").Value = 642102.0428085228
$ws.Range("R7").Value = 6698251.641631705
$ws.Range("S7").Value = 25
$ws.Range("T7").Value = "Uppsala"
$ws.Range("U7").Value = "Tierp"
$ws.Range("V7").Value = "Uppland"
$ws.Range("W7").Value = "Tolfta"
$ws.Range("Y7").NumberFormat = "@"
$ws.Range("Y7").Value = "2023-09-09"
$ws.Range("Z7").NumberFormat = "@"
$ws.Range("Z7").Value = "12:06"
$ws.Range("AA7").NumberFormat = "@"
$ws.Range("AA7").Value = "2023-09-09"
$ws.Range("AB7").NumberFormat = "@"
$ws.Range("AB7").Value = "12:06"
$ws.Range("AD7").Value = $False
$ws.Range("AE7").Value = $False
$ws.Range("AG7").Value = $False
$ws.Range("AW7").Value = "Lotta Lund"
$ws.Range("AX7").Value = "Lotta Lund"

# Row 8
$ws.Range("A8").Value = 111931635
$ws.Range("B8").Value = 90658
$ws.Range("C8").Value = "Ovaliderad"
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 4361
$ws.Range("F8").Value = "Orange taggsvamp"
$ws.Range("G8").Value = "Hydnellum aurantiacum"
$ws.Range("H8").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("P8").Value = "Strömsbergs bruk/spökskogen, Upl"
$ws.Range("Q8").Value = 642301.6473846264
$ws.Range("R8").Value = 6698304.981152782
$ws.Range("S8").Value = 50
$ws.Range("T8").Value = "Uppsala"
$ws.Range("U8").Value = "Tierp"
$ws.Range("V8").Value = "Uppland"
$ws.Range("W8").Value = "Tolfta"
$ws.Range("Y8").NumberFormat = "@"
$ws.Range("Y8").Value = "2023-09-06"
$ws.Range("Z8").NumberFormat = "@"
$ws.Range("Z8").Value = "00:00"
$ws.Range("AA8").NumberFormat = "@"
$ws.Range("AA8").Value = "2023-09-06"
$ws.Range("AB8").NumberFormat = "@"
$ws.Range("AB8").Value = "00:00"
$ws.Range("AD8").Value = $False
$ws.Range("AE8").Value = $False
$ws.Range("AG8").Value = $False
$ws.Range("AH8").Value = "Skogsmark"
$ws.Range("AI8").Value = "Kalkbarrskog"
$ws.Range("AW8").Value = "Lotta Lund"
$ws.Range("AX8").Value = "Lotta Lund"

